$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions - copy the formatting from H1 (bold, centered, bordered)
# onto the two new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 10

$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 9

$ws.Range("I4").Value = 6
$ws.Range("J4").Value = 6
